# Generate Report for Handoff
# Update the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" values
# for the bc7309c9-... file (row 7) across the Overview, zh-cn and de-de sheets,
# reflecting a freshly generated handoff report.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-08-22 06:41:12"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-08-22 06:41:03"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-08-22 06:41:12"
